# Update missing-data summary tables to reflect ln(lac/mann) variable names.
$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("eed_t1-dev_t2")
$ws1.Range("A14").Value = "ln_L_conc_t1"
$ws1.Range("A15").Value = "ln_L_conc_t1"
$ws1.Range("A16").Value = "ln_L_conc_t1"
$ws1.Range("A17").Value = "ln_L_conc_t1"
$ws1.Range("A18").Value = "ln_M_conc_t1"
$ws1.Range("A19").Value = "ln_M_conc_t1"
$ws1.Range("A20").Value = "ln_M_conc_t1"
$ws1.Range("A21").Value = "ln_M_conc_t1"

$ws2 = $wb.Worksheets.Item("eed_t2-dev_t23")
$ws2.Range("A20").Value = "ln_L_conc_t1"
$ws2.Range("A21").Value = "ln_L_conc_t1"
$ws2.Range("A22").Value = "ln_L_conc_t1"
$ws2.Range("A23").Value = "ln_L_conc_t1"
$ws2.Range("A24").Value = "ln_L_conc_t1"
$ws2.Range("A25").Value = "ln_L_conc_t1"
$ws2.Range("A26").Value = "ln_M_conc_t1"
$ws2.Range("A27").Value = "ln_M_conc_t1"
$ws2.Range("A28").Value = "ln_M_conc_t1"
$ws2.Range("A29").Value = "ln_M_conc_t1"
$ws2.Range("A30").Value = "ln_M_conc_t1"
$ws2.Range("A31").Value = "ln_M_conc_t1"

$wb.Save()
